$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 1, 6.201049113329182),
    @(0.127881588408715, 0.04240448674262143, 0.1575252929769615, 0.496779210170732, 0, 0.8245905782990299),
    @(0.127881588408715, 0.3127903958511391, 0.8054896365839992, 0.496779210170732, 0, 1.742940831014585),
    @(0.6753301551942219, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 1, 3.645393585217082),
    @(3.230985683306322, 1.667794583268128, 3.900430680208489, 0.496779210170732, 1, 9.295990156953671),
    @(0.127881588408715, 1.667794583268128, 337.1190423067083, 8.660232485948974, 1, 347.5749509643341),
    @(3.230985683306322, 1.667794583268128, 26.21740644021617, 0.496779210170732, 0, 31.61296591696135),
    @(0.6753301551942219, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 0, 3.645393585217082),
    @(0.127881588408715, 0.04240448674262143, 3.900430680208489, 8.660232485948974, 0, 12.7309492413088),
    @(3.230985683306322, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 1, 5.553084769722144),
    @(3.230985683306322, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 0, 5.553084769722144),
    @(1.459612070389937, 1.667794583268128, 0.8054896365839992, 8.660232485948974, 0, 12.59312877619104),
    @(1.459612070389937, 0.3127903958511391, 0.1575252929769615, 0.496779210170732, 0, 2.42670696938877),
    @(3.230985683306322, 114.8270160096505, 3.900430680208489, 8.660232485948974, 1, 130.6186648591143),
    @(0.6753301551942219, 1.667794583268128, 0.8054896365839992, 645.3272768299601, 0, 648.4758912050064)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]
    for ($j = 0; $j -lt $vals.Length; $j++) {
        $col = 2 + $j
        $ws.Cells.Item($row, $col).Value = $vals[$j]
    }
}
